$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top of the sheet (shifts existing data down by one row)
$ws.Rows.Item(1).Insert()

# Populate the new header row
$ws.Range("A1").Value = "Code"
$ws.Range("B1").Value = "Label"

# Restore the view to the top-left and select B1, matching the target state
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("B1").Select()
